$wb = $excel.ActiveWorkbook

# --- Sheet "Version": update Start Time / End Time / Run Time ---
$wsVersion = $wb.Worksheets.Item("Version")
$wsVersion.Range("C10").Value = "23:27:46 03-05-2017 Sunday Pacific Standard Time"
$wsVersion.Range("C11").Value = "23:28:48 03-05-2017 Sunday Pacific Standard Time"
$wsVersion.Range("C12").Value = "0:1:1 0 days"

# --- Sheet "3_3": update Average/Max/Min/Sdev and Raw data strings ---
$ws33 = $wb.Worksheets.Item("3_3")
$ws33.Range("C3").Value = -0.623
$ws33.Range("D3").Value = -0.597
$ws33.Range("E3").Value = -0.639
$ws33.Range("F3").Value = 0.011
$ws33.Range("H3").Value = "-0.629141,-0.624487,-0.628027,-0.597395,-0.618176,-0.623387,-0.631179,-0.629881,-0.629508,-0.631274,-0.623006,-0.628393,-0.62338,-0.626811,-0.6284,-0.622646,-0.62821,-0.627191,-0.633686,-0.633027,-0.621247,-0.633407,-0.62107,-0.621342,-0.627564,-0.628122,-0.624582,-0.63109,-0.602782,-0.599433,-0.598312,-0.623291,-0.628767,-0.628027,-0.626817,-0.604453,-0.597388,-0.629331,-0.624215,-0.627191,-0.626729,-0.631933,-0.628862,-0.624956,-0.624501,-0.627653,-0.632857,-0.620785,-0.63592,-0.629324,-0.618088,-0.622177,-0.629419,-0.628951,-0.627191,-0.6334,-0.619575,-0.622558,-0.624969,-0.635642,-0.623658,-0.622829,-0.635907,-0.626906,-0.624406,-0.625982,-0.622918,-0.624963,-0.628122,-0.632571,-0.629691,-0.623013,-0.600357,-0.634344,-0.620785,-0.620886,-0.628393,-0.629698,-0.626627,-0.621715,-0.638794,-0.628387,-0.624596,-0.635452,-0.623475,-0.630261,-0.599154,-0.601199,-0.63266,-0.622924,-0.63285,-0.631831,-0.624677,-0.601566,-0.596742,-0.599249,-0.599901,-0.601376,-0.599426,-0.598128"

$ws33.Range("C4").Value = -9.877000000000001
$ws33.Range("D4").Value = -9.856
$ws33.Range("E4").Value = -9.894
$ws33.Range("F4").Value = 0.008999999999999999
$ws33.Range("H4").Value = "-9.869865,-9.871618,-9.877746,-9.86549,-9.867338,-9.874688,-9.887026,-9.881183,-9.875612,-9.858357,-9.87032,-9.882772,-9.888058,-9.869491,-9.871801,-9.88657,-9.883316,-9.8771,-9.869206,-9.874872,-9.891217,-9.882379,-9.866604,-9.865402,-9.866896,-9.88619,-9.87441,-9.865218,-9.86475,-9.878493,-9.883031,-9.87672,-9.868181,-9.881278,-9.892052,-9.875796,-9.875517,-9.861876,-9.88407,-9.886013,-9.876638,-9.864301,-9.874308,-9.886937,-9.878398,-9.87403,-9.862521,-9.87867,-9.882209,-9.876822,-9.86801,-9.868839,-9.880252,-9.889729,-9.885823,-9.862712,-9.870592,-9.879607,-9.890748,-9.869668,-9.860762,-9.872821,-9.889818,-9.884804,-9.863269,-9.873017,-9.881272,-9.891128,-9.877195,-9.864281,-9.874226,-9.887596,-9.887412,-9.865782,-9.865775,-9.883683,-9.894091,-9.872555,-9.866142,-9.868751,-9.883405,-9.890653,-9.882297,-9.868561,-9.872637,-9.891299,-9.882386,-9.874682,-9.867623,-9.881924,-9.885823,-9.890647,-9.856027,-9.875524,-9.882215,-9.885279,-9.879519,-9.870225,-9.874587,-9.890565"

$ws33.Range("C5").Value = -1.074
$ws33.Range("D5").Value = -0.594
$ws33.Range("E5").Value = -10.354
$ws33.Range("F5").Value = 1.776
$ws33.Range("H5").Value = "-0.633319,-0.614283,-5.420498,-0.634059,-0.629603,-0.628393,-0.637306,-0.623937,-0.612524,-0.631267,-0.631084,-0.623203,-0.629242,-0.626817,-0.629508,-1.66741,-0.627646,-5.136143,-0.627007,-0.628482,-0.619113,-0.593679,-0.61912,-10.354458,-0.629698,-0.598964,-0.616512,-0.639731,-0.629324,-0.631457,-0.626729,-0.629698,-4.359852,-0.624311,-0.625798,-0.6309,-0.629236,-0.595533,-0.625343,-0.625526,-0.627273,-0.63535,-0.634901,-0.624589,-0.624963,-0.637673,-0.63499,-1.832675,-0.627931,-0.621709,-0.637028,-0.630255,-0.637857,-0.626824,-0.637402,-0.626349,-0.632938,-0.623563,-0.631736,-0.632021,-0.6259,-9.988316,-0.778889,-0.616879,-0.623753,-0.607225,-0.600819,-0.626729,-0.631185,-0.628584,-0.617993,-0.626729,-0.625153,-0.62376,-0.594698,-0.62162,-0.630812,-1.25569,-0.62357,-0.62997,-0.624399,-0.631274,-0.630356,-0.629983,-0.626349,-0.631749,-0.633971,-0.624786,-0.597293,-0.626817,-0.632673,-0.62857,-0.62088,-0.623658,-0.629521,-0.629976,-0.628407,-0.63249,-0.62626,-10.258656"

# --- Sheet "1_8": update Average/Max/Min/Sdev and Raw data strings ---
$ws18 = $wb.Worksheets.Item("1_8")
$ws18.Range("C3").Value = -0.001
$ws18.Range("D3").Value = 0.008
$ws18.Range("E3").Value = -0.011
$ws18.Range("F3").Value = 0.004
$ws18.Range("H3").Value = "-0.003278,0.003918,0.003628,-0.006349,0.002808,0.001982,0.000323,-0.000793,0.004479,-0.002181,-0.007465,-0.001638,0.001957,-0.004949,0.002259,-0.005214,-0.003593,-0.007459,-0.002729,0.002247,0.001969,-0.000244,-0.003013,-0.002446,0.005577,-0.003001,-0.001083,-0.010795,-0.003858,-0.007724,0.005034,5.2e-05,-0.000799,-0.000522,0.002265,6.4e-05,0.001717,-0.004154,0.000317,0.000607,-0.001077,0.003924,-0.000251,4.5e-05,0.000311,-0.001891,0.001427,0.003671,-0.003821,0.000317,-0.001095,-0.000812,-0.001336,0.000329,0.00253,-0.002729,-0.000251,5.2e-05,-0.000787,-0.006324,-0.004117,-0.001607,-0.00133,-0.000546,0.001667,3.3e-05,-0.000238,-0.002446,0.007809,0.001402,0.00504,-0.003019,-0.001367,0.001149,-0.006059,0.006421,-0.002199,-0.008026,0.003634,0.001162,-0.003025,0.001975,0.001433,-0.001626,-0.001348,-0.000226,-0.003025,0.001717,0.003375,0.005854,-0.002452,0.001692,0.000884,-0.006361,-0.003839,0.004473,-0.001897,0.00837,-0.000805,-0.004665"

$ws18.Range("C4").Value = -0
$ws18.Range("D4").Value = 0.008
$ws18.Range("E4").Value = -0.01
$ws18.Range("F4").Value = 0.003
$ws18.Range("H4").Value = "-0.003568,0.002524,0.0006,-0.00387,-0.000516,0.001667,-0.001064,0.007266,-0.009938,-0.003556,-0.003315,-0.001897,0.001704,0.000329,0.001975,-0.005523,0.000896,-0.005233,0.001994,0.001427,0.003098,-0.002742,-0.006065,0.004738,-0.001915,0.002839,-0.003007,-0.001897,0.000311,0.003092,0.002512,0.001661,0.000317,0.001433,-0.000805,0.005028,0.002808,0.003375,-0.00136,-0.000238,0.00118,0.002247,-0.00329,-0.000812,0.001427,-0.0044,-0.000516,-0.003297,0.001143,0.001149,-0.003278,0.000872,-0.000238,-0.001619,-0.002729,-0.003568,-0.002994,-0.003297,0.00446,-0.002748,-0.003303,0.000323,-0.003864,0.000311,-0.001613,-0.003858,-0.005782,0.00586,-0.003858,0.006421,0.002265,-0.005251,-0.005529,0.003067,-0.000793,-0.003019,-0.002168,-0.001064,-0.003556,0.007778,0.001723,0.001692,0.00475,-0.003019,-0.001909,0.004195,-0.002476,-0.003321,-0.00136,0.003344,-0.001613,0.000582,0.000872,-0.00247,-0.004419,0.001414,0.0006,-0.003574,0.001427,0.000335"

$ws18.Range("C5").Value = -0.001
$ws18.Range("D5").Value = 0.006
$ws18.Range("E5").Value = -0.011
$ws18.Range("F5").Value = 0.003
$ws18.Range("H5").Value = "-0.001619,-0.00133,-0.006602,-0.003013,-0.004678,0.003357,-0.001064,0.004479,-0.004092,0.000341,0.00253,-0.003556,-0.00136,-0.003007,-0.001064,0.003073,0.00253,-0.003303,-0.002717,3.9e-05,-0.003297,0.000317,-0.004961,0.000613,-0.002181,-0.003858,-0.002988,0.005879,0.002802,0.003092,-0.002711,-0.002168,0.002234,-0.000534,0.001975,0.001143,-0.006047,-0.005775,-0.005782,0.001168,0.00311,0.000853,-0.001915,-0.000509,0.003936,-0.000528,0.003634,-0.003025,-0.002156,0.000866,-0.005233,7e-05,-0.003007,-0.00358,-0.002754,-0.004986,-0.003007,-0.011054,-0.004678,-0.000793,0.001162,-0.004141,-0.003315,-0.002735,0.002284,0.00504,-0.005245,-0.001089,0.000335,-0.003833,2.7e-05,-0.00104,0.001131,0.004208,0.00533,-0.00247,-0.005233,-0.006059,0.00171,-0.004394,3.9e-05,0.000884,5.2e-05,0.000866,-0.006343,0.000866,-0.001336,-0.008291,-0.005486,0.001704,-0.004413,-0.001052,-0.007163,-0.000509,0.001414,-0.003568,-0.000207,-0.001632,-0.005202,0.004454"
